$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.867.39"
$ws.Range("E2").Value = "'  +1.73%  "
$ws.Range("D3").Value = "'3.604.89"
$ws.Range("E3").Value = "'  +0.17%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'242.88"
$ws.Range("E5").Value = "'  +3.54%  "
$ws.Range("D6").Value = "'658.65"
$ws.Range("E6").Value = "'  +0.96%  "
$ws.Range("D7").Value = "'1.69"
$ws.Range("E7").Value = "'  +15.76%  "
$ws.Range("E8").Value = "'  +3.16%  "
$ws.Range("E9").Value = "'  +6.48%  "
$ws.Range("D11").Value = "'3.602.40"
$ws.Range("E11").Value = "'  +0.20%  "
$ws.Range("D12").Value = "'43.68"
$ws.Range("E12").Value = "'  +4.04%  "
$ws.Range("E13").Value = "'  +1.73%  "
$ws.Range("E14").Value = "'  +0.89%  "
$ws.Range("D15").Value = "'4.275.04"
$ws.Range("E15").Value = "'  -0.01%  "
$ws.Range("D16").Value = "'96.781.92"
$ws.Range("E17").Value = "'  +2.37%  "
$ws.Range("D18").Value = "'3.609.29"
$ws.Range("E18").Value = "'  +0.44%  "
$ws.Range("D19").Value = "'7.80"
$ws.Range("E19").Value = "'  -0.84%  "
$ws.Range("D20").Value = "'12.71"
$ws.Range("E20").Value = "'  -1.12%  "
$ws.Range("D21").Value = "'18.03"
$ws.Range("E21").Value = "'  +1.02%  "
$ws.Range("D22").Value = "'0.542"
$ws.Range("E22").Value = "'  +14.10%  "
$ws.Range("D23").Value = "'511.80"
$ws.Range("E23").Value = "'  +0.92%  "
$ws.Range("D24").Value = "'3.41"
$ws.Range("E24").Value = "'  -2.35%  "
$ws.Range("D25").Value = "'0.0000202"
$ws.Range("E25").Value = "'  +4.39%  "
$ws.Range("D26").Value = "'6.85"
$ws.Range("E26").Value = "'  +4.42%  "
$ws.Range("D27").Value = "'98.07"
$ws.Range("E27").Value = "'  +3.10%  "
$ws.Range("D28").Value = "'13.06"
$ws.Range("E28").Value = "'  +5.16%  "
$ws.Range("D29").Value = "'3.797.57"
$ws.Range("E29").Value = "'  +0.04%  "
$ws.Range("D30").Value = "'3.04"
$ws.Range("E30").Value = "'  -0.14%  "
$ws.Range("D31").Value = "'0.149"
$ws.Range("E31").Value = "'  +7.67%  "
$ws.Range("D32").Value = "'11.53"
$ws.Range("E32").Value = "'  +3.00%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "'  +0.01%  "
$ws.Range("D34").Value = "'0.185"
$ws.Range("E34").Value = "'  +5.20%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "'  -0.07%  "
$ws.Range("D36").Value = "'31.67"
$ws.Range("E36").Value = "'  -2.15%  "
$ws.Range("D37").Value = "'624.67"
$ws.Range("E37").Value = "'  +12.44%  "
$ws.Range("D38").Value = "'0.571"
$ws.Range("E38").Value = "'  +2.27%  "
$ws.Range("E39").Value = "'  +8.34%  "
$ws.Range("D40").Value = "'1.62"
$ws.Range("E40").Value = "'  +11.02%  "
$ws.Range("D41").Value = "'0.153"
$ws.Range("E41").Value = "'  +1.83%  "
$ws.Range("E42").Value = "'  +0.03%  "
$ws.Range("D43").Value = "'1.90"
$ws.Range("E43").Value = "'  +8.83%  "
$ws.Range("D44").Value = "'0.919"
$ws.Range("E44").Value = "'  +1.22%  "
$ws.Range("D45").Value = "'5.91"
$ws.Range("E45").Value = "'  +4.55%  "
$ws.Range("D46").Value = "'0.0431"
$ws.Range("E46").Value = "'  +4.89%  "
$ws.Range("E47").Value = "'  +2.16%  "
$ws.Range("D48").Value = "'23.68"
$ws.Range("E48").Value = "'  +0.55%  "
$ws.Range("D49").Value = "'33.57"
$ws.Range("E49").Value = "'  -7.28%  "
$ws.Range("D50").Value = "'8.38"
$ws.Range("E50").Value = "'  +4.82%  "
$ws.Range("D51").Value = "'3.52"
$ws.Range("E51").Value = "'  -1.86%  "
